$d = $word.ActiveDocument

# Locate the "Funciones:" paragraph (bold heading before the functions list).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.Trim() -eq "Funciones:") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Funciones:' paragraph"
}

$rng = $target.Range

$bodyXml = @'
<w:p><w:pPr><w:spacing w:after="240" w:lineRule="auto"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Enfoque Inicial:</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> El enfoque planificado consistía en desarrollar simultáneamente tanto el programa de compresión y encriptación como el de descompresión y desencriptación, construyendo progresivamente la complejidad de ambos sistemas. Se inició codificando los componentes básicos necesarios para el funcionamiento, implementando gradualmente funcionalidades adicionales en ambos programas para resolver el desafío. La meta era que, a medida que se </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">incrementara</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> la complejidad, se alcanzaría un punto donde el programa de desencriptado y descompresión coincidiría con los requerimientos del desafío.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="240" w:before="240" w:lineRule="auto"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Cambio de Estrategia:</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Este plan fue abandonado al identificar que los algoritmos que se pretendían implementar para la compresión RLE y LZ78 no coincidían con los utilizados en el dataset proporcionado. Sin embargo, no se </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">anticipaban</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> problemas con los componentes de encriptación/desencriptación.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="240" w:before="240" w:lineRule="auto"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Nuevo Enfoque:</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Se decidió entonces comenzar creando primero las funciones de desencriptación, tomando como referencia específica el archivo README.txt incluido en el dataset del proyecto.</w:t></w:r></w:p><w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Funciones:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:rPr><w:b w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>
'@

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

Write-Output "done"
